$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2024")

# Insert a new row above row 36, shifting existing rows 36-125 down to 37-126.
$ws.Rows.Item(36).Insert()

# Populate the new row 36 with the new September transaction entry.
$ws.Cells.Item(36, 18).Value = "ach indianesign bal axisbank"
$ws.Cells.Item(36, 19).Value = "2024-09-10 13:22:37"
